$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.633.32'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '1.695.98'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '315.27'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = '0.3920'
$ws.Range('E7').Value = '  -0.56%  '
$ws.Range('D8').Value = '0.4038'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('E9').Value = '  -0.56%  '
$ws.Range('D10').Value = '0.9986'
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('D11').Value = '52.96'
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('D12').Value = '0.08830'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').Value = '7.400'
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').Value = '23.55'
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('D15').Value = '8.131'
$ws.Range('E15').Value = '  +6.78%  '
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').Value = '1.699.82'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = '99.44'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').Value = '0.07014'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').Value = '19.70'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').Value = '7.055'
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').Value = '14.59'
$ws.Range('E23').Value = '  +3.65%  '
$ws.Range('D24').Value = '24.630.20'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').Value = '3.118'
$ws.Range('E25').Value = '  +3.29%  '
$ws.Range('D26').Value = '2.345'
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range('D27').Value = '22.61'
$ws.Range('E27').Value = '  +0.92%  '
$ws.Range('D28').Value = '163.06'
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('D29').Value = '8.749'
$ws.Range('E29').Value = '  +17.02%  '
$ws.Range('D30').Value = '135.52'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').Value = '5.142'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').Value = '0.08938'
$ws.Range('D33').Value = '7.572'
$ws.Range('E33').Value = '  +3.06%  '
$ws.Range('D34').Value = '1.064'
$ws.Range('E34').Value = '  -4.03%  '
$ws.Range('D35').Value = '1.960'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '11.02'
$ws.Range('D37').Value = '0.2742'
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').Value = '0.02879'
$ws.Range('E38').Value = '  +3.71%  '
$ws.Range('D39').Value = '14.40'
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('D40').Value = '0.09147'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').Value = '1.458'
$ws.Range('E41').Value = '  -0.60%  '
$ws.Range('D42').Value = '0.7629'
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('D43').Value = '15.83'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('D44').Value = '0.7163'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').Value = '2.573'
$ws.Range('E45').Value = '  +1.22%  '
$ws.Range('D46').Value = '4.204'
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = '1.336'
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('E49').Value = '  -1.08%  '
$ws.Range('D50').Value = '0.07965'
$ws.Range('D51').Value = '90.28'
$ws.Range('E51').Value = '  +1.99%  '
